# Add "Section" and "ID No" columns (H and I) to the Grade1 roster sheet.
# - H1 header "Section", with "-" filled for every data row (2..33), matching
#   the existing "-" placeholder pattern used elsewhere in the sheet.
# - I1 header "ID No", with each student's unique ID number (row 2..33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---------------------------------------------------------------
$ws.Range("H1").Value2 = "Section"
$ws.Range("I1").Value2 = "ID No"

# Match header style (A1:G1 all use style 0).
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# --- Section column (H) : constant "-" placeholder for every data row ------
$ws.Range("H2:H33").Value2 = "-"

# H already existed with style 0 on every row - keep that.
$ws.Range("H2:H33").NumberFormat = "General"

# --- ID No column (I) : unique per-student id -------------------------------
$ids = @{
    2  = "18-0036"
    3  = "16-0137"
    4  = "17-0011"
    5  = "17-0008"
    6  = "17-0044"
    7  = "17-0045"
    8  = "16-0003"
    9  = "16-0031"
    10 = "16-0011"
    11 = "16-0129"
    12 = "16-0030"
    13 = "16-0013"
    14 = "16-0008"
    15 = "16-0159"
    16 = "16-0109"
    17 = "16-0014"
    18 = "16-0054"
    19 = "16-0194"
    20 = "17-0018"
    21 = "17-0010"
    22 = "16-0028"
    23 = "18-0224"
    24 = "16-0052"
    25 = "16-0025"
    26 = "16-0183"
    27 = "17-0049"
    28 = "16-0138"
    29 = "17-0054"
    30 = "17-0079"
    31 = "16-0204"
    32 = "16-0168"
    33 = "16-0154"
}

foreach ($row in $ids.Keys) {
    $idCell = "I$row"
    $ws.Range($idCell).Value2 = $ids[$row]

    # Give the new ID cell the same look & feel as the rest of that row
    # (column B carries the per-row style used throughout the sheet).
    $ws.Range("B$row").Copy()
    $ws.Range($idCell).PasteSpecial(-4122)
    $ws.Range($idCell).Value2 = $ids[$row]
}

# --- Column widths / view state --------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 8.43
$ws.Columns.Item(9).ColumnWidth = 8.43

$excel.ActiveWindow.Zoom = 65
$ws.Range("A1").Select()
$ws.Range("I2:I33").Select()
